$d = $word.ActiveDocument

# Replace all Helvetica font references with Arial across the whole document
$d.Content.Font.Name = $d.Content.Font.Name
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Helvetica"
$find.Replacement.Text = "Arial"
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)

# Update the numeric confidence-interval values that changed
$d.Content.Find.Execute("1.28 (0.66 - 2.38)", $true, $false, $false, $false, $false, $true, 1, $false, "1.28 (0.66 - 2.4)", 2)
$d.Content.Find.Execute("0.17 (0.05 - 0.4)", $true, $false, $false, $false, $false, $true, 1, $false, "0.17 (0.05 - 0.42)", 2)
$d.Content.Find.Execute("0.36 (0.22 - 0.72)", $true, $false, $false, $false, $false, $true, 1, $false, "0.36 (0.22 - 0.76)", 2)
$d.Content.Find.Execute("0.8 (0.48 - 1.76)", $true, $false, $false, $false, $false, $true, 1, $false, "0.8 (0.48 - 1.73)", 2)
